$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "PURCHASER"
$ws.Range("A12").Value = "SALES"

$ws.Range("A12").Select()
